# petty-cashBook-2021.xlsx - 26-Jan-2021 midday update
# Sheet1 ("Buku KAS HARIAN" ledger view) - add new transactions to the daily
# cash book (rows 6-11) and update two existing rows (3-5) with additional
# amounts folded into their formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 3: D3 gains an extra 195,000 ---------------------------------
$ws.Range("D3").Formula = "=45000+195000"

# --- Row 4: D4 gains three more additions/subtraction terms -----------
$ws.Range("D4").Formula = "=49000+37256000+3465000+6240000+1375000+2424000-1215000"

# --- Row 5: C5 gains an extra 28,367,500 -------------------------------
$ws.Range("C5").Formula = "=37256000+28367500"

# --- Row 6 (new): A/P, 1,266,000 out -----------------------------------
$ws.Range("B6").Value = "A/P"
$ws.Range("D6").Formula = "=1266000"

# --- Row 7 (new): PRIVE - andreas, 5,000,000 out (plain value) --------
$ws.Range("B7").Value = "PRIVE - andreas"
$ws.Range("D7").Value = 5000000

# --- Row 8 (new): SALES - cash/retail, in ------------------------------
$ws.Range("B8").Value = "SALES - cash/retail"
$ws.Range("C8").Formula = "=15349475+25076525-28367500"

# --- Row 9 (new): SELISIH - lebih, in ----------------------------------
$ws.Range("B9").Value = "SELISIH - lebih"
$ws.Range("C9").Formula = "=100000-65500"

# --- Row 10 (new): SETOR KE BANK, out ----------------------------------
$ws.Range("B10").Value = "SETOR KE BANK"
$ws.Range("D10").Formula = "=22000000"

# --- Row 11 (new): next day's date entry -------------------------------
$ws.Range("A11").Value = 44222

# --- View state: scroll down a bit and land the selection on C31 ------
$ws.Range("C31").Select()
